$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing data rows (old rows 6 and 7).
# Deleting row 6 twice shifts row 7's content up into row 6, then removes it.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# Rewrite the remaining data rows (2-5) with their updated values.
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F7"
$ws.Range("C2").Value = "F3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.174047333333333
$ws.Range("H2").Value = 3.522142
$ws.Range("I2").Value = 0.9719567830576163
$ws.Range("J2").Value = 0.9719567830576163
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 173.637756
$ws.Range("N2").Value = 520.913268
$ws.Range("O2").Value = 0.9875738673498291
$ws.Range("P2").Value = 0.9875738673498291
$ws.Range("Q2").Value = 203.858944397784
$ws.Range("R2").Value = 1834.730499580056
$ws.Range("S2").Value = 0.959879119141109
$ws.Range("T2").Value = 0.959879119141109

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F7"
$ws.Range("C3").Value = "F3"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.174047333333333
$ws.Range("H3").Value = 3.522142
$ws.Range("I3").Value = 0.9719567830576163
$ws.Range("J3").Value = 0.9719567830576163
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.184794333333333
$ws.Range("N3").Value = 6.554383
$ws.Range("O3").Value = 0.0124261326501708
$ws.Range("P3").Value = 0.0124261326501708
$ws.Range("Q3").Value = 2.565051960931778
$ws.Range("R3").Value = 23.085467648386
$ws.Range("S3").Value = 0.01207766391650723
$ws.Range("T3").Value = 0.01207766391650723

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F7"
$ws.Range("C4").Value = "F3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.033874
$ws.Range("H4").Value = 0.101622
$ws.Range("I4").Value = 0.02804321694238366
$ws.Range("J4").Value = 0.02804321694238367
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 173.637756
$ws.Range("N4").Value = 520.913268
$ws.Range("O4").Value = 0.9875738673498291
$ws.Range("P4").Value = 0.9875738673498291
$ws.Range("Q4").Value = 5.881805346744
$ws.Range("R4").Value = 52.93624812069601
$ws.Range("S4").Value = 0.02769474820872008
$ws.Range("T4").Value = 0.02769474820872009

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "F7"
$ws.Range("C5").Value = "F3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.033874
$ws.Range("H5").Value = 0.101622
$ws.Range("I5").Value = 0.02804321694238366
$ws.Range("J5").Value = 0.02804321694238367
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.184794333333333
$ws.Range("N5").Value = 6.554383
$ws.Range("O5").Value = 0.0124261326501708
$ws.Range("P5").Value = 0.0124261326501708
$ws.Range("Q5").Value = 0.07400772324733333
$ws.Range("R5").Value = 0.6660695092259999
$ws.Range("S5").Value = 0.0003484687336635767
$ws.Range("T5").Value = 0.0003484687336635767
